# This script inserts one new data row into the sheet at row 312, shifting
# the existing rows 312..441 down to 313..442 (which matches a full
# shift of all subsequent rows by one position, as described by the diff),
# and fills the newly inserted row 312 with its own data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 312; this shifts rows 312:441 down to 313:442,
# carrying along their values and styles.
$ws.Rows.Item(312).Insert()

# Populate the newly inserted row 312 with the new record's data.
# Columns A,B,C,E,F,G,H,I,N,O,Q,R mirror the same "template" values used
# throughout this block of rows (Terminal La Palmera de La Serena / Espinaca).
$ws.Cells.Item(312, 1).Value = 8
$ws.Cells.Item(312, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(312, 3).Value = "Coquimbo"
$ws.Cells.Item(312, 4).Value = 45119
$ws.Cells.Item(312, 5).Value = 4
$ws.Cells.Item(312, 6).Value = 100112012
$ws.Cells.Item(312, 7).Value = "Espinaca"
$ws.Cells.Item(312, 8).Value = "Sin especificar"
$ws.Cells.Item(312, 9).Value = "Primera"
$ws.Cells.Item(312, 10).Value = 400
$ws.Cells.Item(312, 11).Value = 500
$ws.Cells.Item(312, 12).Value = 600
$ws.Cells.Item(312, 13).Value = 550
$ws.Cells.Item(312, 14).Value = "$/atado 300 a 500 gramos"
$ws.Cells.Item(312, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(312, 16).Value = 1100
$ws.Cells.Item(312, 17).Value = 0.5
$ws.Cells.Item(312, 18).Value = "Hortaliza"
